$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "no coupling" test vector previously carried stale o1/o2/o3 sample
# values (a constant 4204391, a row index 1..7, and a lingering output
# comment style). Reset the F3:H9 block so the expected outputs are all
# zero, and make sure the three output columns share one consistent
# center/bottom alignment instead of the old mismatched styles.
$rng = $ws.Range("F3:H9")
$rng.HorizontalAlignment = -4108   # xlCenter
$rng.VerticalAlignment = -4107     # xlBottom

For ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 6).Value = 0.0   # F -> o1
    $ws.Cells.Item($r, 7).Value = 0.0   # G -> o2
    $ws.Cells.Item($r, 8).Value = 0.0   # H -> o3
}
